{"js": "// The document contains a line reading \"Version:10.1.2\" (the Nessus tool\n// version, in paragraph \"...3. Nessus: / [tab][tab]Version:10.1.2\").\n// The edit bumps the version to 10.6.1. In the canonical OOXML this is\n// expressed by splitting the final run (\":10.1.2\", rFonts hint=\"eastAsia\" +\n// bCs) into two runs with identical formatting: \":10.\" and \"6.1\".\n//\n// We reproduce that exact run split: locate the \"1.2\" tail of the run,\n// momentarily toggle a character formatting property (bold) to force Word\n// to break it into its own run, replace its text with \"6.1\", then restore\n// the original (non-bold) formatting on the newly created run so the final\n// run properties match the source (rFonts hint=\"eastAsia\", bCs only).\n\nconst body = context.document.body;\n\n// Find the unique \"Version:10.1.2\" run in the document.\nconst matches = body.search(\":10.1.2\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \":10.1.2\" in the document.');\n}\n\nconst fullMatch = matches.items[0];\nconst para = fullMatch.paragraphs.getFirst();\n\n// Narrow down to the \"1.2\" part (the run split happens right before it,\n// turning \":10.1.2\" into \":10.\" + \"1.2\").\nconst tailMatches = fullMatch.search(\"1.2\", { matchCase: true });\ntailMatches.load(\"items\");\nawait context.sync();\n\nconst tail = tailMatches.items[0];\n\n// Force a run boundary at the split point by toggling formatting.\ntail.font.bold = true;\ntail.insertText(\"6.1\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-acquire the just-inserted text (scoped to the same paragraph, so it\n// stays unambiguous even if \"6.1\" happens to appear elsewhere) and restore\n// its formatting so the final run keeps the same (non-bold) character\n// formatting as the rest of the run it was split from.\nconst newTailMatches = para.search(\"6.1\", { matchCase: true });\nnewTailMatches.load(\"items\");\nawait context.sync();\n\nconst newTail = newTailMatches.items[0];\nnewTail.font.bold = false;\nawait context.sync();\n", "ps1": "# The document contains a line reading \"Version:10.1.2\" (the Nessus tool\n# version, in the paragraph that starts with \"3. Nessus:\" followed by\n# \"[tab][tab]Version:10.1.2\"). The fix bumps the version to 10.6.1.\n#\n# In the canonical OOXML this shows up as the trailing run (\":10.1.2\",\n# formatted with rFonts hint=\"eastAsia\" + bCs) being split into two runs\n# with identical formatting: \":10.\" and \"6.1\".\n#\n# We reproduce that exact run split: find the unique \":10.1.2\" text, narrow\n# to its \"1.2\" tail, momentarily toggle Bold to force Word to break that\n# tail into its own run, set its text to \"6.1\", then restore Bold to its\n# original (off) state so the final run's formatting matches the source.\n\n$d = $word.ActiveDocument\n\n# Locate the unique \":10.1.2\" run in the document.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \":10.1.2\"\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Narrow to the \"1.2\" tail within that match (this is where the run\n    # split happens: \":10.\" stays put, \"1.2\" becomes its own run that we\n    # rewrite).\n    $tail = $rng.Duplicate\n    $tail.Find.ClearFormatting()\n    $tail.Find.Text = \"1.2\"\n    $tail.Find.MatchCase = $true\n    $tail.Find.Execute() | Out-Null\n\n    $tail.Font.Bold = 1\n    $tail.Text = \"6.1\"\n    $tail.Font.Bold = 0\n}\n"}
